# Apply cryptocurrency price/volume updates per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.976.15"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "3.370.28"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  -0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "555.48"
$c.NumberFormat = "General"
$ws.Range("E5").Value = "  -0.07%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "174.59"
$c.NumberFormat = "General"
$ws.Range("E6").Value = "  -0.77%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.631"
$c.NumberFormat = "General"
$ws.Range("E7").Value = "  +1.88%  "
$ws.Range("D8").Value = "3.360.24"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  +5.36%  "
$ws.Range("E11").Value = "  +1.14%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "53.79"
$c.NumberFormat = "General"
$ws.Range("E12").Value = "  -1.62%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000280"
$c.NumberFormat = "General"
$ws.Range("E13").Value = "  +2.41%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "9.17"
$c.NumberFormat = "General"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").Value = "3.903.70"
$ws.Range("E15").Value = "  +0.42%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "18.28"
$c.NumberFormat = "General"
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").Value = "3.366.57"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").Value = "65.034.73"
$ws.Range("E20").Value = "  +0.91%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.NumberFormat = "General"
$ws.Range("E21").Value = "  +1.35%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "454.21"
$c.NumberFormat = "General"
$ws.Range("E22").Value = "  +0.08%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "4.92"
$c.NumberFormat = "General"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("E24").Value = "  -0.48%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "87.51"
$c.NumberFormat = "General"
$ws.Range("E25").Value = "  +2.46%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "13.87"
$c.NumberFormat = "General"
$ws.Range("E26").Value = "  +3.41%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.88"
$c.NumberFormat = "General"
$ws.Range("E27").Value = "  +1.19%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "10.74"
$c.NumberFormat = "General"
$ws.Range("E28").Value = "  -0.51%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "8.68"
$c.NumberFormat = "General"
$ws.Range("E29").Value = "  -1.14%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "31.20"
$c.NumberFormat = "General"
$ws.Range("E30").Value = "  +4.10%  "
$ws.Range("E31").Value = "  -1.36%  "
$ws.Range("E32").Value = "  -0.50%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "62.85"
$c.NumberFormat = "General"
$ws.Range("E33").Value = "  +7.43%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "576.39"
$c.NumberFormat = "General"
$ws.Range("E34").Value = "  -0.44%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.107"
$c.NumberFormat = "General"
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  +5.12%  "
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("E39").Value = "  +0.09%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.372"
$c.NumberFormat = "General"
$ws.Range("E40").Value = "  +0.98%  "
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("D42").Value = "3.077.82"
$ws.Range("E42").Value = "  -0.77%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.0415"
$c.NumberFormat = "General"
$ws.Range("E43").Value = "  +1.34%  "
$ws.Range("E44").Value = "  -1.16%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "3.20"
$c.NumberFormat = "General"
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.46"
$c.NumberFormat = "General"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.134"
$c.NumberFormat = "General"
$ws.Range("E47").Value = "  +2.37%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "142.31"
$c.NumberFormat = "General"
$ws.Range("E48").Value = "  +5.06%  "
$ws.Range("E49").Value = "  +0.02%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.55"
$c.NumberFormat = "General"
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("E51").Value = "  -0.72%  "
